# Added support for apiKeys
# Fill in "Done" markers (columns C:F) for the /events, /inbox and
# /inbox/unread rows on the Methods worksheet, matching the pattern
# already used for the other method rows on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

$rows = @(19, 85, 86)
foreach ($r in $rows) {
    $ws.Range("C$r`:F$r").Value = "Done"
}

# Update the view state to match where the author left the selection.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 70
$ws.Range("C86").Select()
